# Nerfs to Merchant Class Specialties
# - Reduce the scaling of skill_bonus (column J) for rows 3-15
# - Remove the now-unused cost (column W) values for rows 2-15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Affixes")

# New skill_bonus (column J) values, keyed by row number
$jValues = [ordered]@{
    3  = 0.07
    4  = 0.12
    5  = 0.18
    6  = 0.24
    7  = 0.29
    8  = 0.35
    9  = 0.41
    10 = 0.47
    11 = 0.52
    12 = 0.58
    13 = 0.64
    14 = 0.69
    15 = 0.75
}

foreach ($row in $jValues.Keys) {
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}

# Clear the cost (column W) cells entirely for rows 2-15
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 23).ClearContents()
}
